$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add Q3 (empty cell, same style as P3) ---
$ws.Range("P3").Copy()
$ws.Range("Q3").PasteSpecial(-4122)

# --- Add Q4 = 2020 (same style as P4, but vertical alignment top) ---
$ws.Range("P4").Copy()
$q4 = $ws.Range("Q4")
$q4.PasteSpecial(-4122)
$q4.Value = 2020
$q4.VerticalAlignment = -4160

# --- Q5 = 1.1 (same style as existing, s=13 already, just set value) ---
$ws.Range("Q5").Value = 1.1000000000000001

# --- Q6 = 7 (same style as A3/style 9 cells, but with 0.0 number format) ---
$ws.Range("A3").Copy()
$q6 = $ws.Range("Q6")
$q6.PasteSpecial(-4122)
$q6.Value = 7
$q6.NumberFormat = "0.0"

# --- Update selection to J22 ---
$ws.Range("J22").Select()

$excel.CutCopyMode = 0
